$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 86
$ws.Range("A$row").Value = "Golang Developer"
$ws.Range("B$row").Value = "https://www.dice.com/job-detail/9a35a2c6-39b2-4363-8018-5a763af2592a"
$ws.Range("C$row").Value = "West Chester, Pennsylvania"
$ws.Range("D$row").Value = "Third Party, Contract"
$ws.Range("E$row").Value = "Depends on Experience"
$ws.Range("F$row").Value = "InfiCare Technologies"
